$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 275-279 are rotated: HBA (previously row 279) moves to the top (row 275),
# and MEL/NOU/PER/SYD each shift down by one row.
$ws.Range("A275").Value = "HBA"
$ws.Range("B275").Value = "Hobart, Australia"
$ws.Range("C275").Value = -42.883209
$ws.Range("D275").Value = 147.331665
$ws.Range("E275").Value = "AU"
$ws.Range("F275").Value = "Oceania"
$ws.Range("G275").Value = "Hobart"

$ws.Range("A276").Value = "MEL"
$ws.Range("B276").Value = "Melbourne, VIC, Australia"
$ws.Range("C276").Value = -37.6733016968
$ws.Range("D276").Value = 144.843002319
$ws.Range("E276").Value = "AU"
$ws.Range("F276").Value = "Oceania"
$ws.Range("G276").Value = "Melbourne"

$ws.Range("A277").Value = "NOU"
$ws.Range("B277").Value = "Noumea, New Caledonia"
$ws.Range("C277").Value = -22.0146007538
$ws.Range("D277").Value = 166.212997436
$ws.Range("E277").Value = "NC"
$ws.Range("F277").Value = "Oceania"
$ws.Range("G277").Value = "Noumea"

$ws.Range("A278").Value = "PER"
$ws.Range("B278").Value = "Perth, WA, Australia"
$ws.Range("C278").Value = -31.9402999878
$ws.Range("D278").Value = 115.967002869
$ws.Range("E278").Value = "AU"
$ws.Range("F278").Value = "Oceania"
$ws.Range("G278").Value = "Perth"

$ws.Range("A279").Value = "SYD"
$ws.Range("B279").Value = "Sydney, NSW, Australia"
$ws.Range("C279").Value = -33.9460983276
$ws.Range("D279").Value = 151.177001953
$ws.Range("E279").Value = "AU"
$ws.Range("F279").Value = "Oceania"
$ws.Range("G279").Value = "Sydney"

# Row 284 (Christchurch) gets updated lat/lon coordinates.
$ws.Range("C284").Value = -43.4893989563
$ws.Range("D284").Value = 172.5319976807
